# Generate Report for Handoff
#
# The underlying data for file "91bbbb67-038f-4402-9f8c-f2af61374c0d.md" has
# progressed (new handoff) while "ea39ce26-1598-49c4-bb04-a987ab57c00b.md"
# remains in its previously handed-back state. Because each sheet lists rows
# by status/date ordering, the two files swap row order (ea39ce26 now first,
# 91bbbb67 now second) and 91bbbb67's status/date columns are refreshed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "ea39ce26-1598-49c4-bb04-a987ab57c00b.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "2016-03-24 10:19:21"

$ws.Range("A3").Value = "91bbbb67-038f-4402-9f8c-f2af61374c0d.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-03-24 10:21:18"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "ea39ce26-1598-49c4-bb04-a987ab57c00b.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "ea39ce26-1598-49c4-bb04-a987ab57c00b.88569866682b099058049a19a9b0edb5f426a9c1.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-24 10:19:13"
$ws.Range("F2").Value = "ea39ce26-1598-49c4-bb04-a987ab57c00b.md"
$ws.Range("G2").Value = "ea39ce26-1598-49c4-bb04-a987ab57c00b.88569866682b099058049a19a9b0edb5f426a9c1.zh-cn.xlf"
$ws.Range("H2").Value = "2016-03-24 10:19:59"
$ws.Range("J2").Value = "Include"

$ws.Range("A3").Value = "91bbbb67-038f-4402-9f8c-f2af61374c0d.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "91bbbb67-038f-4402-9f8c-f2af61374c0d.c31ca651c21b93c2e8ff5d2f0b01ba6f1462c3e1.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-24 10:21:09"
$ws.Range("F3").Value = "91bbbb67-038f-4402-9f8c-f2af61374c0d.md"
$ws.Range("G3").Value = "91bbbb67-038f-4402-9f8c-f2af61374c0d.c31ca651c21b93c2e8ff5d2f0b01ba6f1462c3e1.zh-cn.xlf"
$ws.Range("H3").Value = "2016-03-24 10:19:59"
$ws.Range("J3").Value = "Include"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "ea39ce26-1598-49c4-bb04-a987ab57c00b.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "ea39ce26-1598-49c4-bb04-a987ab57c00b.88569866682b099058049a19a9b0edb5f426a9c1.de-de.xlf"
$ws.Range("E2").Value = "2016-03-24 10:19:21"
$ws.Range("F2").Value = "ea39ce26-1598-49c4-bb04-a987ab57c00b.md"
$ws.Range("G2").Value = "ea39ce26-1598-49c4-bb04-a987ab57c00b.88569866682b099058049a19a9b0edb5f426a9c1.de-de.xlf"
$ws.Range("H2").Value = "2016-03-24 10:20:16"
$ws.Range("J2").Value = "Include"

$ws.Range("A3").Value = "91bbbb67-038f-4402-9f8c-f2af61374c0d.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "91bbbb67-038f-4402-9f8c-f2af61374c0d.c31ca651c21b93c2e8ff5d2f0b01ba6f1462c3e1.de-de.xlf"
$ws.Range("E3").Value = "2016-03-24 10:21:18"
$ws.Range("F3").Value = "91bbbb67-038f-4402-9f8c-f2af61374c0d.md"
$ws.Range("G3").Value = "91bbbb67-038f-4402-9f8c-f2af61374c0d.c31ca651c21b93c2e8ff5d2f0b01ba6f1462c3e1.de-de.xlf"
$ws.Range("H3").Value = "2016-03-24 10:20:16"
$ws.Range("J3").Value = "Include"
